# Word COM-interop script reproducing the target edit:
#
#   "...le Tout Puissant et ses Puissants..."
#        -> "...le Tout-Puissant et ses Puissants..."
#
# with the space between "Tout" and "Puissant" turned into a hyphen,
# the text split across three runs ("...le Tout" / "-" / "Puissant...")
# and the document's "_GoBack" bookmark (which Word drops at the last
# edited location) left sitting right between the new "-" run and the
# following "Puissant" run.  Word only ever keeps a single "_GoBack"
# bookmark, so re-adding it here also removes the old one that used to
# sit after "vis-à-vis de" later in the document.

$d = $word.ActiveDocument

# Locate the first (and relevant) occurrence of "le Tout Puissant".
$found = $d.Content
$ok = $found.Find.Execute("le Tout Puissant", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) {
    throw "Could not find 'le Tout Puissant' in the document"
}
$matchStart = $found.Start

# Force a run boundary right after "Tout" (before the space) so the
# following edits don't get silently re-merged into a single run.
$splitPoint = $d.Range($matchStart + 7, $matchStart + 7)
$d.Bookmarks.Add("_TempSplit", $splitPoint)

# Replace the space between "Tout" and "Puissant" with a hyphen; this
# becomes its own run thanks to the boundary created above.
$spaceRange = $d.Range($matchStart + 7, $matchStart + 8)
$spaceRange.Text = "-"

# Remove the temporary helper bookmark again.
$d.Bookmarks("_TempSplit").Delete()

# Drop the "_GoBack" bookmark right between the new hyphen run and
# "Puissant" - this also relocates/replaces whatever "_GoBack" bookmark
# already existed elsewhere in the document.
$goBackPoint = $d.Range($matchStart + 8, $matchStart + 8)
$d.Bookmarks.Add("_GoBack", $goBackPoint)
